{"js": "// 1) Merge the split runs that make up:\n//    Use the \"UI Trait\" (which is to specify a component name and\n// into a single run of text (the quotes around \"UI Trait\" are curly\n// quotes: \\u201C and \\u201D).\nconst openQuote = \"\\u201C\";\nconst closeQuote = \"\\u201D\";\nconst mergedText =\n  \"Use the \" + openQuote + \"UI Trait\" + closeQuote +\n  \" (which is to specify a component name and \";\n\nconst useTheResults = context.document.body.search(\n  \"Use the \" + openQuote + \"UI Trait\" + closeQuote +\n    \" (which is to specify a component name and \",\n  { matchCase: true }\n);\nuseTheResults.load(\"text\");\nawait context.sync();\n\nif (useTheResults.items.length > 0) {\n  useTheResults.items[0].insertText(mergedText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Fix the typo \"For and component\" -> \"For a component\".\nconst typoResults = context.document.body.search(\n  \"For and component, we need to know\",\n  { matchCase: true }\n);\ntypoResults.load(\"text\");\nawait context.sync();\n\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\n    \"For a component, we need to know\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Merge the split runs that make up:\n#    Use the \"UI Trait\" (which is to specify a component name and\n# into a single run of text (quotes are curly: U+201C / U+201D).\n$openQuote = [char]0x201C\n$closeQuote = [char]0x201D\n$mergedText = \"Use the \" + $openQuote + \"UI Trait\" + $closeQuote + \" (which is to specify a component name and \"\n\n$find1 = $d.Content.Find\n$find1.Text = $mergedText\n$find1.Replacement.Text = $mergedText\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 2) Fix the typo \"For and component\" -> \"For a component\".\n$find2 = $d.Content.Find\n$find2.Text = \"For and component, we need to know\"\n$find2.Replacement.Text = \"For a component, we need to know\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
